$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("F56").Value = 0
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("F61").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("F67").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("F78").Value = 0
$ws.Range("F79").Value = 0
$ws.Range("F80").Value = 0
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 0
$ws.Range("F83").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("F90").Value = 0
$ws.Range("F91").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("D113").Value = 0
$ws.Range("E113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("F152").Value = 0
$ws.Range("F153").Value = 0
$ws.Range("F154").Value = 0
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("F158").Value = 0
$ws.Range("F159").Value = 0
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("F163").Value = 0
$ws.Range("F164").Value = 0
$ws.Range("F165").Value = 0
$ws.Range("F166").Value = 0
$ws.Range("F167").Value = 0
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("F170").Value = 0
$ws.Range("F171").Value = 0
$ws.Range("F172").Value = 0
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("F174").Value = 0
$ws.Range("F175").Value = 0
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("F178").Value = 0
$ws.Range("F179").Value = 0
$ws.Range("F180").Value = 0
$ws.Range("F181").Value = 0
$ws.Range("F182").Value = 0
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("F186").Value = 0
$ws.Range("F187").Value = 0
$ws.Range("F188").Value = 0
$ws.Range("F189").Value = 0
$ws.Range("F190").Value = 0
$ws.Range("F191").Value = 0
$ws.Range("F192").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("F194").Value = 0
$ws.Range("F195").Value = 0
$ws.Range("F196").Value = 0
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("F199").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("F204").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("F213").Value = 0
